# Populate "import.xlsx" style test workbook with a third data table and a
# Sheet3 label, then leave the workbook focused on Sheet3 (matching the
# updated test fixture: Sheet2 gains the same header/data/footer table that
# Sheet1 has, Sheet3 gains a "Sheet 3" label, and the active tab moves to
# Sheet3).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet2: duplicate Sheet1's header + data rows + merged footer row ----
$headers = @("a1", "b1", "c1", "d1")
$row2    = @("a2", "b2", "c2", "d2")
$row3    = @("a3", "b3", "c3", "d3")
$row4    = @("a4", "b4", "c4", "d4")

$cols = @("A", "B", "C", "D")

for ($i = 0; $i -lt 4; $i++) {
    $ws2.Range($cols[$i] + "1").Value = $headers[$i]
    $ws2.Range($cols[$i] + "2").Value = $row2[$i]
    $ws2.Range($cols[$i] + "3").Value = $row3[$i]
    $ws2.Range($cols[$i] + "4").Value = $row4[$i]
}

$ws2.Range("A5").Value = "5ad"
$ws2.Range("A5:D5").Merge()
$ws2.Range("A5:D5").HorizontalAlignment = -4108
$ws2.Range("A5:D5").VerticalAlignment = -4108

# --- Sheet3: single label cell --------------------------------------------
$ws3.Range("A1").Value = "Sheet 3"

# --- Selection / active-sheet bookkeeping ---------------------------------
$ws1.Range("A1").Select()
$ws2.Range("A1").Select()
$ws3.Range("A2").Select()

$ws3.Activate()
